$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$hl = $ws.Range("F5").Hyperlinks.Item(1)
$hl.TextToDisplay = "ttest@example.com"
